# "add Png and modify makeObject"
# Remap the values in the two stamped "makeObject" tile-grids (A6:P12) on
# Sheet1: 4 -> 3, 2 -> 6, 0 -> 7. Then leave the selection where the user
# last clicked while making the edit (J11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$map = @{ 4 = 3; 2 = 6; 0 = 7 }

for ($r = 6; $r -le 12; $r++) {
    for ($c = 1; $c -le 16; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $old = [int]$cell.Value2
        if ($map.ContainsKey($old)) {
            $cell.Value = $map[$old]
        }
    }
}

$ws.Range("J11").Select()
